# "termino de importacion de orientacion espacial"
# The question bank's `correct_answer` column (F) was still holding the
# raw option index (1/2/3/4) from the earlier partial import; finish the
# import by switching it over to the option letter (A/B/C/D) that the
# rest of the sheet (answer_a..answer_d / feedback_text) already expects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2:F6 -> correct_answer, numeric option index => option letter
$ws.Range("F2").Value = "A"
$ws.Range("F3").Value = "B"
$ws.Range("F4").Value = "C"
$ws.Range("F5").Value = "D"
$ws.Range("F6").Value = "A"

# Finishing touches that come from re-confirming the cell formatting for
# the whole imported range (locked/protected cells, wrapped feedback text).
$ws.Range("A1:H6").Locked = $true
$ws.Range("G2:G4").WrapText = $true

# Minor re-anchoring of the row-marker images (they get fitted to the new
# row content) - nudge each to its recalculated size.
$EMU = 12700.0

$ws.Shapes.Item(1).Width  = 166680 / $EMU
$ws.Shapes.Item(1).Height = 170280 / $EMU

$ws.Shapes.Item(2).Width  = 126360 / $EMU
$ws.Shapes.Item(2).Height = 132840 / $EMU

$ws.Shapes.Item(3).Width  = 133920 / $EMU
$ws.Shapes.Item(3).Height = 140040 / $EMU

$ws.Shapes.Item(4).Width  = 141120 / $EMU
$ws.Shapes.Item(4).Height = 141840 / $EMU

$ws.Shapes.Item(5).Width  = 138600 / $EMU
$ws.Shapes.Item(5).Height = 140040 / $EMU

$ws.Shapes.Item(6).Width  = 169200 / $EMU
$ws.Shapes.Item(6).Height = 172800 / $EMU

$ws.Shapes.Item(7).Width  = 147960 / $EMU
$ws.Shapes.Item(7).Height = 155520 / $EMU

$ws.Shapes.Item(8).Width  = 127440 / $EMU
$ws.Shapes.Item(8).Height = 133200 / $EMU

# Import finished - leave the cursor parked below the data, back at col A.
$ws.Range("A9").Select()
